$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 214.53
$ws.Range("I15").Value = 214.53
$ws.Range("K15").Value = 643.59
$ws.Range("M15").Value = -474.59

# Hunk 1: ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 343.41666
$ws.Range("I33").Value = 256.83334
$ws.Range("K33").Value = 256.83334
$ws.Range("M33").Value = -27.83334000000002

# Hunk 2: ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1356.2273
$ws.Range("I98").Value = 1402.1765
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 1402.1765
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = 95.82349999999997
$ws.Range("N98").Value = -4196

# Hunk 3: ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 30306650
$ws.Range("I116").Value = 3843.5715
$ws.Range("J116").Value = 83336560
$ws.Range("K116").Value = 3843.5715
$ws.Range("L116").Value = 83336560
$ws.Range("M116").Value = -401.5715
$ws.Range("N116").Value = -83343444

# Hunk 4: ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1356.2273
$ws.Range("I122").Value = 1402.1765
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4206.529500000001
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1756.529500000001
$ws.Range("N122").Value = -8500

# Hunk 5: ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1625.0212
$ws.Range("I132").Value = 1944.0714
$ws.Range("J132").Value = 1154.8422
$ws.Range("K132").Value = 5832.2142
$ws.Range("L132").Value = 3464.5266
$ws.Range("M132").Value = -3302.2142
$ws.Range("N132").Value = -8524.526600000001

# Hunk 6: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1784.1945
$ws.Range("I137").Value = 1080.6595
$ws.Range("J137").Value = 3106.84
$ws.Range("K137").Value = 3241.9785
$ws.Range("L137").Value = 9320.52
$ws.Range("M137").Value = -691.9785000000002
$ws.Range("N137").Value = -14420.52

# Hunk 7: ARM row 14
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 506
$ws.Range("I14").Value = 506
$ws.Range("K14").Value = 506
$ws.Range("M14").Value = -331

# Hunk 8: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6947.2583
$ws.Range("I32").Value = 5886
$ws.Range("J32").Value = 16852.334
$ws.Range("K32").Value = 5886
$ws.Range("L32").Value = 16852.334
$ws.Range("M32").Value = -5599
$ws.Range("N32").Value = -17426.334

# Hunk 9: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5391.52
$ws.Range("I61").Value = 4725
$ws.Range("K61").Value = 4725
$ws.Range("M61").Value = -4513

# Hunk 10: ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1550.5
$ws.Range("I63").Value = 1334
$ws.Range("J63").Value = 2200
$ws.Range("K63").Value = 1334
$ws.Range("L63").Value = 2200
$ws.Range("M63").Value = -648
$ws.Range("N63").Value = -3572

# Hunk 11: ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1550.5
$ws.Range("I66").Value = 1334
$ws.Range("J66").Value = 2200
$ws.Range("K66").Value = 6670
$ws.Range("L66").Value = 11000
$ws.Range("M66").Value = -3238
$ws.Range("N66").Value = -17864

# Hunk 12: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4881.081
$ws.Range("I74").Value = 3047.0833
$ws.Range("J74").Value = 8266.923000000001
$ws.Range("K74").Value = 3047.0833
$ws.Range("L74").Value = 8266.923000000001
$ws.Range("M74").Value = -2173.0833
$ws.Range("N74").Value = -10014.923

# Hunk 13: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4881.081
$ws.Range("I77").Value = 3047.0833
$ws.Range("J77").Value = 8266.923000000001
$ws.Range("K77").Value = 15235.4165
$ws.Range("L77").Value = 41334.61500000001
$ws.Range("M77").Value = -10867.4165
$ws.Range("N77").Value = -50070.61500000001

# Hunk 14: ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2556.5
$ws.Range("I102").Value = 2372.375
$ws.Range("K102").Value = 2372.375
$ws.Range("M102").Value = -750.375

# Hunk 15: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5391.52
$ws.Range("I136").Value = 4725
$ws.Range("K136").Value = 14175
$ws.Range("M136").Value = -11625

# Hunk 16: BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11496480
$ws.Range("I86").Value = 11496480
$ws.Range("K86").Value = 11496480
$ws.Range("M86").Value = -11495357

# Hunk 17: BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 11496480
$ws.Range("I89").Value = 11496480
$ws.Range("K89").Value = 57482400
$ws.Range("M89").Value = -57476784

# Hunk 18: BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1973.1666
$ws.Range("I99").Value = 1725
$ws.Range("J99").Value = 2256.7856
$ws.Range("K99").Value = 1725
$ws.Range("L99").Value = 2256.7856
$ws.Range("M99").Value = -227
$ws.Range("N99").Value = -5252.7856

# Hunk 19: BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4720.773
$ws.Range("I105").Value = 4382.4814
$ws.Range("J105").Value = 5258.0586
$ws.Range("K105").Value = 4382.4814
$ws.Range("L105").Value = 5258.0586
$ws.Range("M105").Value = -2635.4814
$ws.Range("N105").Value = -8752.0586

# Hunk 20: BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1930.5238
$ws.Range("I107").Value = 1790.0625
$ws.Range("K107").Value = 1790.0625
$ws.Range("M107").Value = 129.9375

# Hunk 21: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3023.7874
$ws.Range("I134").Value = 2930.8484
$ws.Range("J134").Value = 3242.8572
$ws.Range("K134").Value = 8792.5452
$ws.Range("L134").Value = 9728.571599999999
$ws.Range("M134").Value = -6257.5452
$ws.Range("N134").Value = -14798.5716

# Hunk 22: CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2015.25
$ws.Range("I31").Value = 1643.921
$ws.Range("J31").Value = 2656.6365
$ws.Range("K31").Value = 1643.921
$ws.Range("L31").Value = 2656.6365
$ws.Range("M31").Value = -1348.921
$ws.Range("N31").Value = -3246.6365

# Hunk 23: CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2015.25
$ws.Range("I34").Value = 1643.921
$ws.Range("J34").Value = 2656.6365
$ws.Range("K34").Value = 1643.921
$ws.Range("L34").Value = 2656.6365
$ws.Range("M34").Value = -1441.921
$ws.Range("N34").Value = -3060.6365

# Hunk 24: CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1360.2142
$ws.Range("I94").Value = 1655.75
$ws.Range("J94").Value = 1242
$ws.Range("K94").Value = 1655.75
$ws.Range("L94").Value = 1242
$ws.Range("M94").Value = -1204.75
$ws.Range("N94").Value = -2144

# Hunk 25: CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 358.2
$ws.Range("I107").Value = 377.81818
$ws.Range("J107").Value = 304.25
$ws.Range("K107").Value = 377.81818
$ws.Range("L107").Value = 304.25
$ws.Range("M107").Value = 1542.18182
$ws.Range("N107").Value = -4144.25

# Hunk 26: CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 12643.772
$ws.Range("I122").Value = 3628.1428
$ws.Range("J122").Value = 28421.125
$ws.Range("K122").Value = 10884.4284
$ws.Range("L122").Value = 85263.375
$ws.Range("M122").Value = -8434.428400000001
$ws.Range("N122").Value = -90163.375

# Hunk 27: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2667.2222
$ws.Range("I132").Value = 2059
$ws.Range("J132").Value = 3701.2
$ws.Range("K132").Value = 6177
$ws.Range("L132").Value = 11103.6
$ws.Range("M132").Value = -3647
$ws.Range("N132").Value = -16163.6

# Hunk 28: CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 640.4583
$ws.Range("I122").Value = 518.38464
$ws.Range("J122").Value = 784.7273
$ws.Range("K122").Value = 4665.46176
$ws.Range("L122").Value = 7062.545700000001
$ws.Range("M122").Value = -2215.46176
$ws.Range("N122").Value = -11962.5457

# Hunk 29: CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17998.072
$ws.Range("I131").Value = 504.42307
$ws.Range("J131").Value = 33682.035
$ws.Range("K131").Value = 1513.26921
$ws.Range("L131").Value = 101046.105
$ws.Range("M131").Value = 3526.73079
$ws.Range("N131").Value = -111126.105

# Hunk 30: GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3149.9355
$ws.Range("I102").Value = 3621.8333
$ws.Range("J102").Value = 2496.5386
$ws.Range("K102").Value = 3621.8333
$ws.Range("L102").Value = 2496.5386
$ws.Range("M102").Value = -1999.8333
$ws.Range("N102").Value = -5740.5386

# Hunk 31: GSM row 104
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 49500
$ws.Range("J104").Value = 49500
$ws.Range("L104").Value = 49500
$ws.Range("N104").Value = -56488

# Hunk 32: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5256.9287
$ws.Range("I122").Value = 6549.7
$ws.Range("J122").Value = 2025
$ws.Range("K122").Value = 19649.1
$ws.Range("L122").Value = 6075
$ws.Range("M122").Value = -17199.1
$ws.Range("N122").Value = -10975

# Hunk 33: LTW row 50
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 13940.429
$ws.Range("J50").Value = 13940.429
$ws.Range("L50").Value = 13940.429
$ws.Range("N50").Value = -15214.429

# Hunk 34: LTW row 54
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Hunk 35: LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 26981.889
$ws.Range("I61").Value = 35638.668
$ws.Range("J61").Value = 9668.333000000001
$ws.Range("K61").Value = 35638.668
$ws.Range("L61").Value = 9668.333000000001
$ws.Range("M61").Value = -35436.668
$ws.Range("N61").Value = -10072.333

# Hunk 36: LTW row 63
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 32542.5
$ws.Range("J63").Value = 32542.5
$ws.Range("L63").Value = 32542.5
$ws.Range("N63").Value = -34040.5

# Hunk 37: LTW row 66
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 32542.5
$ws.Range("J66").Value = 32542.5
$ws.Range("L66").Value = 97627.5
$ws.Range("N66").Value = -105115.5

# Hunk 38: LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 26981.889
$ws.Range("I113").Value = 35638.668
$ws.Range("J113").Value = 9668.333000000001
$ws.Range("K113").Value = 35638.668
$ws.Range("L113").Value = 9668.333000000001
$ws.Range("M113").Value = -33468.668
$ws.Range("N113").Value = -14008.333

# Hunk 39: LTW row 137
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 51799.168
$ws.Range("J137").Value = 51799.168
$ws.Range("L137").Value = 51799.168
$ws.Range("N137").Value = -61999.168

# Hunk 40: WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2924.9285
$ws.Range("I122").Value = 1605.421
$ws.Range("K122").Value = 4816.263
$ws.Range("M122").Value = -2366.263

# Hunk 41: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1053
$ws.Range("I132").Value = 998.30554
$ws.Range("J132").Value = 1217.0834
$ws.Range("K132").Value = 2994.91662
$ws.Range("L132").Value = 3651.2502
$ws.Range("M132").Value = -464.91662
$ws.Range("N132").Value = -8711.2502

# Hunk 42: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3283.605
$ws.Range("I136").Value = 1575.5625
$ws.Range("J136").Value = 5768.0303
$ws.Range("K136").Value = 4726.6875
$ws.Range("L136").Value = 17304.0909
$ws.Range("M136").Value = -2176.6875
$ws.Range("N136").Value = -22404.0909
